# "Jail checks refactored completely - a few to add later."
# Appends 20 new data rows (1583-1602) to the bottom of the case-data table
# on Sheet1, directly after the existing last row (1582).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each inner array is: RowNumber, A, B, C, D, E, F, G, H, I, J, K
# A $null entry for G means "no value in that column at all" (cell omitted).
$rows = @(
     ,@(1583, '21TRD09386', 'Hemmeter', 'DUS Ucm', '4510.111', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0', '5', 'None')
    ,@(1584, '21TRD09386', 'Hemmeter', 'Tail Lights-rear License Plate', '4513.05', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None')
    ,@(1585, '21TRD09386', 'Hemmeter', 'DUS Ucm', '4510.111', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0', '5', 'None')
    ,@(1586, '21TRD09386', 'Hemmeter', 'Tail Lights-rear License Plate', '4513.05', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None')
    ,@(1587, '21TRD09386', 'Hemmeter', 'DUS Ucm', '4510.111', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0', '5', 'None')
    ,@(1588, '21TRD09386', 'Hemmeter', 'Tail Lights-rear License Plate', '4513.05', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None')
    ,@(1589, '21TRD09386', 'Hemmeter', 'DUS Ucm', '4510.111', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0', '5', 'None')
    ,@(1590, '21TRD09386', 'Hemmeter', 'Tail Lights-rear License Plate', '4513.05', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None')
    ,@(1591, '21TRD09386', 'Hemmeter', 'DUS Ucm', '4510.111', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0', '5', 'None')
    ,@(1592, '21TRD09386', 'Hemmeter', 'Tail Lights-rear License Plate', '4513.05', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None')
    ,@(1593, '21TRD09386', 'Bunner', 'DUS Ucm', '4510.111', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0', '5', 'None')
    ,@(1594, '21TRD09386', 'Bunner', 'Tail Lights-rear License Plate', '4513.05', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None')
    ,@(1595, '21TRD09386', 'Bunner', 'DUS Ucm', '4510.111', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0', '5', 'None')
    ,@(1596, '21TRD09386', 'Bunner', 'Tail Lights-rear License Plate', '4513.05', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None')
    ,@(1597, '21TRC08418', 'Bunner', 'Driving In Marked Lanes', '4511.33', 'MM', 'Dismissed', $null, ' ', ' ', ' ', ' ')
    ,@(1598, '21TRC08418', 'Bunner', 'Turn And Stop Signals', 'No Data', 'MM', 'Dismissed', $null, ' ', ' ', ' ', ' ')
    ,@(1599, '21TRC08418', 'Bunner', 'OVI Alcohol / Drugs 1st', '4511.19A1A*', 'M1', 'Guilty', 'Guilty', '$ 375', '$ 0', '180', '177')
    ,@(1600, '21TRC08418', 'Bunner', 'Driving In Marked Lanes', '4511.33', 'MM', 'Dismissed', '__EMPTY__', ' ', ' ', ' ', ' ')
    ,@(1601, '21TRC08418', 'Bunner', 'Turn And Stop Signals', 'No Data', 'MM', 'Dismissed', '__EMPTY__', ' ', ' ', ' ', ' ')
    ,@(1602, '21TRC08418', 'Bunner', 'OVI Alcohol / Drugs 1st', '4511.19A1A*', 'M1', 'Guilty', 'Guilty', '$ 375', '$ 0', '180', '177')
)

foreach ($row in $rows) {
    $r = $row[0]

    # Columns A-F and H-K are written as plain text in the source data (case
    # numbers, statute codes like "4510.111", counts like "5"/"180", and
    # "$ 0" style money strings all must stay text, not get coerced to
    # numbers/dates).
    $ws.Range("A$r`:F$r").NumberFormat = "@"
    $ws.Range("H$r`:K$r").NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]

    $g = $row[7]
    if ($g -eq $null) {
        # leave column G entirely blank/untouched for this row (no cell
        # written at all)
    } elseif ($g -eq '__EMPTY__') {
        # G participates in the row's span but carries no value -- still
        # touch it (format only) so a blank cell entry is emitted
        $ws.Cells.Item($r, 7).NumberFormat = "@"
    } else {
        $ws.Cells.Item($r, 7).NumberFormat = "@"
        $ws.Cells.Item($r, 7).Value = $g
    }

    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
}
